# Update "想去人数" (interested-count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 1779
$ws1.Range("F16").Value = 2956
$ws1.Range("F20").Value = 531
$ws1.Range("F26").Value = 1219
$ws1.Range("F28").Value = 1408
$ws1.Range("F41").Value = 122

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 1779
$ws4.Range("F16").Value = 2956
$ws4.Range("F20").Value = 531
$ws4.Range("F27").Value = 1219
$ws4.Range("F29").Value = 1408
$ws4.Range("F43").Value = 122
